$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast-error table values (naive component forecaster bug fix)
$data = @{
    2  = @(0.1560082246840916, 1.334666554170477, 7.673042805661654, 2.770025777075306, 2.793148460838132, 51)
    3  = @(0.07663322156689945, 1.299834471760893, 7.386217108190144, 2.7177595751262, 2.744260142950399, 50)
    4  = @(0.1568606132576995, 1.29112821977527, 7.117402227092668, 2.667845990137487, 2.690829533573063, 49)
    5  = @(0.0935746678326199, 1.365793217520896, 7.57656095386177, 2.752555349827097, 2.78007587260341, 48)
    6  = @(0.1424348649439729, 1.381738402566604, 8.024261160333362, 2.832712685807257, 2.859715526557411, 47)
    7  = @(0.1066018923057719, 1.373320158004856, 7.768055166292034, 2.787123098517903, 2.815859037392981, 46)
    8  = @(0.1578286088130434, 1.342716586736134, 7.844528394684448, 2.800808525173481, 2.827956351813759, 45)
    9  = @(0.08991622227243078, 1.431264513475095, 8.13824447033133, 2.852760850532573, 2.884308066253469, 44)
    10 = @(0.157804923366503, 1.414187432328667, 8.224904023511264, 2.867909347157135, 2.897454014080703, 43)
    11 = @(0.06053529240133091, 1.461861534950335, 8.447652183675055, 2.906484506009804, 2.941077792346796, 42)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - ME
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - MAE
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - MSE
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - RMSE
    $ws.Cells.Item($row, 6).Value = $vals[4]   # F - SE
    $ws.Cells.Item($row, 7).Value = $vals[5]   # G - N
}
